$d = $word.ActiveDocument

# --- Collapse multi-run text into single runs (Title, Author, Abstract) ---
$d.Content.Find.Execute("Factsheet: Beta distribution", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Factsheet: Beta distribution", 2) | Out-Null

$d.Content.Find.Execute("Michelle Arnetta and Tom Coleman", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Michelle Arnetta and Tom Coleman", 2) | Out-Null

$d.Content.Find.Execute("A factsheet about the beta distribution.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "A factsheet about the beta distribution.", 2) | Out-Null

# --- Fix m:d delimiter property order (sepChr before endChr) in the math equations ---
# Paragraph indices (by current Word paragraph numbering) that contain affected m:d elements.
$mathFixes = @{
  8 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="BodyText" /></w:pPr><w:r><w:rPr><w:b /><w:bCs /></w:rPr><w:t xml:space="preserve">Notation:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>X</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∼</m:t></m:r><m:r><m:rPr><m:nor /><m:sty m:val="p" /></m:rPr><m:t>Beta</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>α</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>,</m:t></m:r><m:r><m:t>β</m:t></m:r></m:e></m:d></m:oMath></w:p>
'@
  17 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="Compact" /><w:jc w:val="left" /></w:pPr><m:oMath><m:r><m:rPr><m:sty m:val="p" /><m:scr m:val="double-struck" /></m:rPr><m:t>E</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>X</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>α</m:t></m:r></m:num><m:den><m:r><m:t>α</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>β</m:t></m:r></m:den></m:f></m:oMath></w:p>
'@
  21 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="Compact" /><w:jc w:val="left" /></w:pPr><m:oMath><m:r><m:rPr><m:sty m:val="p" /><m:scr m:val="double-struck" /></m:rPr><m:t>V</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>X</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>α</m:t></m:r><m:r><m:t>β</m:t></m:r></m:num><m:den><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>α</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>β</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>α</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>β</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:den></m:f></m:oMath></w:p>
'@
  25 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="Compact" /><w:jc w:val="left" /></w:pPr><m:oMath><m:r><m:rPr><m:sty m:val="p" /><m:scr m:val="double-struck" /></m:rPr><m:t>P</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>X</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>α</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:sup></m:sSup><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d></m:e><m:sup><m:r><m:t>β</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:sup></m:sSup></m:num><m:den><m:r><m:rPr><m:nor /><m:sty m:val="p" /></m:rPr><m:t>B</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>α</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>,</m:t></m:r><m:r><m:t>β</m:t></m:r></m:e></m:d></m:den></m:f></m:oMath></w:p>
'@
  26 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="Compact" /><w:jc w:val="left" /></w:pPr><m:oMath><m:r><m:rPr><m:nor /><m:sty m:val="p" /></m:rPr><m:t>B</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>,</m:t></m:r><m:r><m:t>y</m:t></m:r></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">is the beta function</w:t></w:r></w:p>
'@
  29 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="Compact" /><w:jc w:val="left" /></w:pPr><m:oMath><m:r><m:rPr><m:sty m:val="p" /><m:scr m:val="double-struck" /></m:rPr><m:t>P</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>X</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>≤</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:sSub><m:e><m:r><m:t>I</m:t></m:r></m:e><m:sub><m:r><m:t>x</m:t></m:r></m:sub></m:sSub><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>α</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>,</m:t></m:r><m:r><m:t>β</m:t></m:r></m:e></m:d></m:oMath></w:p>
'@
  30 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="Compact" /><w:jc w:val="left" /></w:pPr><m:oMath><m:sSub><m:e><m:r><m:t>I</m:t></m:r></m:e><m:sub><m:r><m:t>x</m:t></m:r></m:sub></m:sSub><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>,</m:t></m:r><m:r><m:t>b</m:t></m:r></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">is the regularized incomplete beta function</w:t></w:r></w:p>
'@
  35 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="FirstParagraph" /></w:pPr><w:r><w:t xml:space="preserve">Then the distribution of the probabilities of a customer purchasing from Cantor’s Confectionery can be expressed as</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>X</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>∼</m:t></m:r><m:r><m:rPr><m:nor /><m:sty m:val="p" /></m:rPr><m:t>Beta</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>7</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>,</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:d></m:oMath><w:r><w:t xml:space="preserve">, meaning the first shape parameter is 7 and the second shape parameter is 5.</w:t></w:r></w:p>
'@
}

foreach ($idx in ($mathFixes.Keys | Sort-Object)) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertXML($mathFixes[$idx])
}
